# Fixes the sound specification sheet of Scarlet (Scarlet_Wald_M block, rows 21-32):
# Corrects timestamps / details and blanks out three rows (28, 31, 32) that are
# no longer needed, while keeping the cell formatting of the columns that still
# carry a style (C/D "Combat"/"Kampfgeräusch" columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: was "< 30% HP" w/ 00:00.0-00:13.3, becomes "Angriff" w/ 00:22.0-00:22.6
$ws.Range("A21").Value = "00:22.0"
$ws.Range("A21").NumberFormat = "h:mm"
$ws.Range("B21").Value = "00:22.6"
$ws.Range("E21").Value = "Angriff"

# Row 22
$ws.Range("A22").Value = "00:22.7"
$ws.Range("B22").Value = "00:23.4"

# Row 23
$ws.Range("A23").Value = "00:23.6"
$ws.Range("B23").Value = "00:24.0"

# Row 24
$ws.Range("A24").Value = "00:24.2"
$ws.Range("B24").Value = "00:24.7"

# Row 25: was "Angriff", becomes "Einstecken Leicht"
$ws.Range("A25").Value = "00:28.8"
$ws.Range("A25").NumberFormat = "h:mm"
$ws.Range("B25").Value = "00:29.6"
$ws.Range("E25").Value = "Einstecken Leicht"

# Row 26
$ws.Range("A26").Value = "00:29.9"
$ws.Range("B26").Value = "00:30.5"

# Row 27
$ws.Range("A27").Value = "00:30.8"
$ws.Range("B27").Value = "00:31.7"

# Row 28 is removed entirely (A, B, E stay empty; C/D keep their format but lose
# their value).
$ws.Range("A28").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()

# Row 29: was "Einstecken Leicht", becomes "Einstecken Schwer"
$ws.Range("A29").Value = "00:37.3"
$ws.Range("B29").Value = "00:38.5"
$ws.Range("E29").Value = "Einstecken Schwer"

# Row 30
$ws.Range("A30").Value = "00:40.5"
$ws.Range("B30").Value = "00:41.9"

# Row 31 is removed entirely.
$ws.Range("A31").ClearContents()
$ws.Range("B31").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()

# Row 32 is removed entirely.
$ws.Range("A32").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("E32").ClearContents()

# Update the saved view/selection state to match the author's final cursor position.
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E45").Select()
